$d = $word.ActiveDocument

# The last paragraph ("Matodova promena ") loses its trailing space, and a
# long run of new glossary paragraphs is appended after it (plus two blank
# paragraphs). We replace everything from that final paragraph onward in one
# shot via Range.InsertXML with a raw OOXML fragment: Word always folds the
# very last paragraph mark of an end-of-document insertion into the
# document's existing trailing paragraph mark, so the fragment below starts
# by re-stating "Matodova promena" (without the trailing space) and ends
# with an empty <w:p/> that becomes the document's final empty paragraph.
# Building it this way keeps run/paragraph boundaries, the xml:space="preserve"
# runs, the split "S" / "idici promena" run and the lastRenderedPageBreak
# marker byte-for-byte faithful to the target revision.
$last = $d.Paragraphs.Last.Range
$last.Collapse(0)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Matodová proměná</w:t></w:r></w:p><w:p><w:r><w:t>Programovaní:</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t>Operátor – logické funkce znaménka plus, rovná se, minus, atd., něco co provádí operandy</w:t></w:r></w:p><w:p><w:r><w:t>Binární operátor – má dva operandy &lt;=, &gt;= atd</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Unární operand – jeden </w:t></w:r></w:p><w:p><w:r><w:t>Ternární operand  - tři, je jediný zkrácený if = podmínka ? true : false;</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Příkaz – instrukce</w:t></w:r></w:p><w:p><w:r><w:t>Příkaz nepodmíněný – udělá se vždy</w:t></w:r></w:p><w:p><w:r><w:t>Příkaz podmíněný – udělá se pouze pokud platí podmínka</w:t></w:r></w:p><w:p><w:r><w:t>Neúplná podmínka – chybí else</w:t></w:r></w:p><w:p><w:r><w:t>Úplná podmínka – má i else</w:t></w:r></w:p><w:p><w:r><w:t>Složená podmínka – má if, else if a else</w:t></w:r></w:p><w:p><w:r><w:t>Přepínač (switch) – vyhodnotí se podmínka a je jich libovolné množství = case pokud case odpovídá podmínce tak se provede, většinou se do case dává break ale pokud ho tam nedáme tak se vyhodnocují další casy, Default se provede i kdyby neplatil žádný case</w:t></w:r></w:p><w:p><w:r><w:t>Cyklus (loop) – while, for, do-while, cyklus se opakuje dokud platí podmínka, dá se ukončit return nebo break</w:t></w:r></w:p><w:p><w:r><w:t>S</w:t></w:r><w:r><w:t>ídící proměná – proměná na, která se provádí iterace a vyhodnocuje se podmínka</w:t></w:r></w:p><w:p><w:r><w:t>Iterace - posloupnost</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$last.InsertXML($xml)
